$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.275.29'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.14%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.867.41'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.38%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.0000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7224'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.07%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '240.99'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.43%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07832'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.66%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3090'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.64%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.25'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.21%  '

$ws.Range('E11').Value = '  +1.27%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.867.27'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.50%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7218'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.27%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.236'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.58%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '90.77'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.79%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.305.32'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.01%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.855'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.40%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '243.81'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.59%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007813'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.66%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.21'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.44%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.104.22'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.48%  '

$ws.Range('E22').Value = '  -0.03%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.958'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.97%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.00%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1596'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +11.46%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '161.90'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.12%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.962'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.19%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.23'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.03%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.346'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.51%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.495'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.68%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.402'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.67%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.101'
$ws.Range('D32').Style = 'Normal'

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05197'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.79%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.933'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.01%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.186'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.59%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7288'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.65%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.683'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.21%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01856'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.65%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.700'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.36%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.173.95'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.01%  '

$ws.Range('E41').Value = '  -0.89%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.109'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.62%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '72.48'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.58%  '

$ws.Range('E44').Value = '  -0.04%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.84'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.65%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5280'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.31%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.004.29'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.48%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.781'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.21%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.888'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.76%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.289'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.91%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4268'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.89%  '
